# Correct typos & update offloading fig
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix typo: missing period after the closing quote of "BCHW" in the
# Description cell for the "layout" parameter.
$ws.Range("C2").Value = "Desired data layout format, accepted values are ""HWC"", ""CHW"", ""BHWC"", ""BCHW"".`nDefaults to ""HWC""."

# Update the selected/active cell shown when the sheet is viewed.
$ws.Range("C3").Select()
